$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.819.36'
$ws.Cells.Item(2, 5).Value = '  -1.07%  '
$ws.Cells.Item(3, 4).Value = '3.506.42'
$ws.Cells.Item(3, 5).Value = '  +0.00%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '602.92'
$ws.Cells.Item(5, 5).Value = '  -1.13%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.91'
$ws.Cells.Item(6, 5).Value = '  -3.52%  '
$ws.Cells.Item(7, 4).Value = '3.504.72'
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 5).Value = '  -1.63%  '
$ws.Cells.Item(10, 5).Value = '  -1.30%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '7.80'
$ws.Cells.Item(11, 5).Value = '  +2.28%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.422'
$ws.Cells.Item(12, 5).Value = '  -2.35%  '
$ws.Cells.Item(13, 5).Value = '  -1.50%  '
$ws.Cells.Item(14, 4).Value = '4.098.91'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '31.44'
$ws.Cells.Item(15, 5).Value = '  -3.80%  '
$ws.Cells.Item(16, 4).Value = '3.503.32'
$ws.Cells.Item(16, 5).Value = '  +0.03%  '
$ws.Cells.Item(17, 4).Value = '66.820.01'
$ws.Cells.Item(17, 5).Value = '  -0.87%  '
$ws.Cells.Item(18, 5).Value = '  -0.83%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '10.61'
$ws.Cells.Item(19, 5).Value = '  +7.04%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.38'
$ws.Cells.Item(20, 5).Value = '  -3.18%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '15.39'
$ws.Cells.Item(21, 5).Value = '  -1.15%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '434.07'
$ws.Cells.Item(22, 5).Value = '  -3.18%  '
$ws.Cells.Item(23, 5).Value = '  -3.63%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '79.85'
$ws.Cells.Item(24, 5).Value = '  +2.13%  '
$ws.Cells.Item(25, 4).Value = '3.639.99'
$ws.Cells.Item(25, 5).Value = '  -0.17%  '
$ws.Cells.Item(26, 5).Value = '  +0.03%  '
$ws.Cells.Item(27, 2).Value = 'PEPE'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.0000119'
$ws.Cells.Item(27, 5).Value = '  -5.43%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.83'
$ws.Cells.Item(28, 5).Value = '  -2.77%  '
$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.25'
$ws.Cells.Item(29, 5).Value = '  -6.80%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.50'
$ws.Cells.Item(30, 5).Value = '  -0.73%  '
$ws.Cells.Item(31, 2).Value = 'Fetch.AI'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.59'
$ws.Cells.Item(31, 5).Value = '  -4.88%  '
$ws.Cells.Item(32, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.00'
$ws.Cells.Item(32, 5).Value = '  +0.03%  '
$ws.Cells.Item(33, 2).Value = 'Kaspa'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.166'
$ws.Cells.Item(33, 5).Value = '  -1.97%  '
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '25.29'
$ws.Cells.Item(34, 5).Value = '  -1.73%  '
$ws.Cells.Item(35, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(35, 4).Value = '3.496.65'
$ws.Cells.Item(35, 5).Value = '  -0.12%  '
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.79'
$ws.Cells.Item(36, 5).Value = '  -4.67%  '
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.87'
$ws.Cells.Item(37, 5).Value = '  -5.22%  '
$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '8.01'
$ws.Cells.Item(38, 5).Value = '  -0.64%  '
$ws.Cells.Item(39, 2).Value = 'USDe'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.00'
$ws.Cells.Item(39, 5).Value = '  +0.00%  '
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.999'
$ws.Cells.Item(40, 5).Value = '  -0.08%  '
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0888'
$ws.Cells.Item(41, 5).Value = '  -0.97%  '
$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '169.40'
$ws.Cells.Item(42, 5).Value = '  -2.42%  '
$ws.Cells.Item(43, 2).Value = 'Stacks'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.08'
$ws.Cells.Item(43, 5).Value = '  -9.95%  '
$ws.Cells.Item(44, 2).Value = 'Filecoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.42'
$ws.Cells.Item(44, 5).Value = '  -1.45%  '
$ws.Cells.Item(45, 2).Value = 'Mantle'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.897'
$ws.Cells.Item(45, 5).Value = '  +1.48%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '28.88'
$ws.Cells.Item(46, 5).Value = '  -4.07%  '
$ws.Cells.Item(47, 2).Value = 'OKB'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '45.66'
$ws.Cells.Item(47, 5).Value = '  -2.48%  '
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.31'
$ws.Cells.Item(48, 5).Value = '  +0.04%  '
$ws.Cells.Item(49, 2).Value = 'Cosmos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.46'
$ws.Cells.Item(49, 5).Value = '  -2.83%  '
$ws.Cells.Item(50, 2).Value = 'dogwifhat'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.40'
$ws.Cells.Item(50, 5).Value = '  -4.93%  '
$ws.Cells.Item(51, 2).Value = 'SuiNetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.980'
$ws.Cells.Item(51, 5).Value = '  -1.55%  '
